$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 168.25
$ws.Range("I12").Value = 168.25
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 168.25
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 1.75
$ws.Range("N12").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1376.6
$ws.Range("I20").Value = 1376.6
$ws.Range("K20").Value = 1376.6
$ws.Range("M20").Value = -1146.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 12697.833
$ws.Range("I21").Value = 9546.75
$ws.Range("K21").Value = 9546.75
$ws.Range("M21").Value = -9078.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 12697.833
$ws.Range("I23").Value = 9546.75
$ws.Range("K23").Value = 9546.75
$ws.Range("M23").Value = -9312.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 177.8
$ws.Range("I28").Value = 163.14285
$ws.Range("K28").Value = 163.14285
$ws.Range("M28").Value = 321.85715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 347.25
$ws.Range("I29").Value = 255.6
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 766.8
$ws.Range("L29").Value = 1500
$ws.Range("M29").Value = -485.8
$ws.Range("N29").Value = -2062

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 1376.6
$ws.Range("I35").Value = 1376.6
$ws.Range("K35").Value = 1376.6
$ws.Range("M35").Value = -997.5999999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 172.16667
$ws.Range("I38").Value = 172.16667
$ws.Range("K38").Value = 516.50001
$ws.Range("M38").Value = -144.50001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1886.8889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32496

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -102480

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1570.7142
$ws.Range("I127").Value = 1332.5
$ws.Range("K127").Value = 3997.5
$ws.Range("M127").Value = 962.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2625.2
$ws.Range("I2").Value = 778.25
$ws.Range("J2").Value = 10013
$ws.Range("K2").Value = 778.25
$ws.Range("L2").Value = 10013
$ws.Range("M2").Value = -665.25
$ws.Range("N2").Value = -10239

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5444.1377
$ws.Range("I32").Value = 1367.591
$ws.Range("K32").Value = 1367.591
$ws.Range("M32").Value = -1080.591

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1420.8334
$ws.Range("I61").Value = 1420.8334
$ws.Range("K61").Value = 1420.8334
$ws.Range("M61").Value = -1208.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1476.25
$ws.Range("I74").Value = 1399.4286
$ws.Range("J74").Value = 2014
$ws.Range("K74").Value = 1399.4286
$ws.Range("L74").Value = 2014
$ws.Range("M74").Value = -525.4286
$ws.Range("N74").Value = -3762

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1476.25
$ws.Range("I77").Value = 1399.4286
$ws.Range("J77").Value = 2014
$ws.Range("K77").Value = 6997.143
$ws.Range("L77").Value = 10070
$ws.Range("M77").Value = -2629.143
$ws.Range("N77").Value = -18806

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 308.875
$ws.Range("I110").Value = 351.14285
$ws.Range("J110").Value = 13
$ws.Range("K110").Value = 351.14285
$ws.Range("L110").Value = 13
$ws.Range("M110").Value = 1693.85715
$ws.Range("N110").Value = -4103

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2625.2
$ws.Range("I116").Value = 778.25
$ws.Range("J116").Value = 10013
$ws.Range("K116").Value = 778.25
$ws.Range("L116").Value = 10013
$ws.Range("M116").Value = 1515.75
$ws.Range("N116").Value = -14601

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3993.8333
$ws.Range("I132").Value = 3789.8
$ws.Range("K132").Value = 11369.4
$ws.Range("M132").Value = -8839.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1420.8334
$ws.Range("I136").Value = 1420.8334
$ws.Range("K136").Value = 4262.5002
$ws.Range("M136").Value = -1712.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2625.2
$ws.Range("I3").Value = 778.25
$ws.Range("J3").Value = 10013
$ws.Range("K3").Value = 778.25
$ws.Range("L3").Value = 10013
$ws.Range("M3").Value = -664.25
$ws.Range("N3").Value = -10241

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7945.875
$ws.Range("I134").Value = 7242.2666
$ws.Range("K134").Value = 21726.7998
$ws.Range("M134").Value = -19191.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 527.5
$ws.Range("I6").Value = 527.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 527.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -414.5
$ws.Range("N6").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1400
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1400
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2550.1538
$ws.Range("I132").Value = 2550.1538
$ws.Range("K132").Value = 7650.4614
$ws.Range("M132").Value = -5120.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 50000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 50000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 50000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -55060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4063.2222
$ws.Range("I134").Value = 3070.875
$ws.Range("J134").Value = 4481.0527
$ws.Range("K134").Value = 9212.625
$ws.Range("L134").Value = 13443.1581
$ws.Range("M134").Value = -6677.625
$ws.Range("N134").Value = -18513.1581

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 2999.75
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2999.75
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8999.25
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9337.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 899.6
$ws.Range("I34").Value = 371
$ws.Range("J34").Value = 1692.5
$ws.Range("K34").Value = 1113
$ws.Range("L34").Value = 5077.5
$ws.Range("M34").Value = -1029
$ws.Range("N34").Value = -5245.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 619.75
$ws.Range("I55").Value = 619.75
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 1859.25
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -1682.25
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 143752.58
$ws.Range("I134").Value = 143752.58
$ws.Range("K134").Value = 431257.74
$ws.Range("M134").Value = -426187.74

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2355.625
$ws.Range("I140").Value = 1302.5
$ws.Range("K140").Value = 3907.5
$ws.Range("M140").Value = 1272.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 2410.75
$ws.Range("I141").Value = 2410.75
$ws.Range("K141").Value = 7232.25
$ws.Range("M141").Value = -2052.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 21000
$ws.Range("J18").Value = 21000
$ws.Range("L18").Value = 21000
$ws.Range("N18").Value = -21586

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6359.7646
$ws.Range("J43").Value = 16000
$ws.Range("L43").Value = 16000
$ws.Range("N43").Value = -16302

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 8333.333000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16250
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 74000.2
$ws.Range("J123").Value = 74000.2
$ws.Range("L123").Value = 74000.2
$ws.Range("N123").Value = -78900.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1996
$ws.Range("I132").Value = 1996
$ws.Range("K132").Value = 5988
$ws.Range("M132").Value = -3458

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4598.6924
$ws.Range("J46").Value = 2228.125
$ws.Range("L46").Value = 2228.125
$ws.Range("N46").Value = -2604.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2100
$ws.Range("J82").Value = 2500
$ws.Range("L82").Value = 2500
$ws.Range("N82").Value = -3222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2100
$ws.Range("J85").Value = 2500
$ws.Range("L85").Value = 2500
$ws.Range("N85").Value = -4996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8884
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3343331.2
$ws.Range("J15").Value = 14997.5
$ws.Range("L15").Value = 14997.5
$ws.Range("N15").Value = -15573.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5999.6665
$ws.Range("I81").Value = 4000
$ws.Range("K81").Value = 8000
$ws.Range("M81").Value = -6939

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5999.6665
$ws.Range("I84").Value = 4000
$ws.Range("K84").Value = 40000
$ws.Range("M84").Value = -34696

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1596.8
$ws.Range("I132").Value = 1596.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4790.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2260.4
$ws.Range("N132").ClearContents()
